# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2210
#   *_new  -> *_FV2304
# Then wrap the data range in an Excel Table (ListObject) and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null) {
            $s = [string]$val
            if ($s.EndsWith("_old")) {
                $cell.Value = $s.Substring(0, $s.Length - 4) + "_FV2210"
            } elseif ($s.EndsWith("_new")) {
                $cell.Value = $s.Substring(0, $s.Length - 4) + "_FV2304"
            }
        }
    }
}

# Turn the header + data range into a native Excel Table.
$tableRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split below row 1, freeze top pane).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
